$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 19: Speed Limit A3
# New row 20: Speed Limit A61
# Values are written in an order that reproduces the original authoring
# order of the shared-string table (A19, B19, F19, F20, B20, A20).
$ws.Range("A19").Value = "speedLimitA3"
$ws.Range("B19").Value = "Speed Limit A3"
$ws.Range("C19").Value = 2015
$ws.Range("F19").Value = "Potential Speed Limit of 130 kmph on German Autobahn A3"

$ws.Range("F20").Value = "Potential Speed Limit of 130 kmph on German Autobahn A61"
$ws.Range("C20").Value = 2015
$ws.Range("B20").Value = "Speed Limit A61"
$ws.Range("A20").Value = "speedLimitA61"

# Match row heights used for the wrapped description text.
$ws.Rows.Item(19).RowHeight = 30
$ws.Rows.Item(20).RowHeight = 30

# Update the view so the new rows are visible and A19 is selected,
# mirroring the author's final sheet view state.
$ws.Range("A19").Select()
$excel.ActiveWindow.ScrollRow = 16
